$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.720.44"
$ws.Range("E2").Value = "  +1.82%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.578.31"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.13"
$ws.Range("E5").Value = "  +4.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.17"
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.574.68"
$ws.Range("E7").Value = "  +1.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.617"
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("E10").Value = "  +3.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.45"
$ws.Range("E11").Value = "  +9.71%  "
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.94"
$ws.Range("E13").Value = "  -1.55%  "
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.155.27"
$ws.Range("E16").Value = "  -1.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "617.60"
$ws.Range("E17").Value = "  -2.17%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.586.27"
$ws.Range("E18").Value = "  +1.05%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.794.26"
$ws.Range("E19").Value = "  +1.90%  "
$ws.Range("E20").Value = "  -2.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.48"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.890"
$ws.Range("E22").Value = "  -0.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.39"
$ws.Range("E23").Value = "  -16.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.07"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.75"
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.82"
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.65"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.62"
$ws.Range("E29").Value = "  +1.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.16"
$ws.Range("E30").Value = "  -2.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.49"
$ws.Range("E31").Value = "  -1.28%  "
$ws.Range("E32").Value = "  -4.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.04"
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("E34").Value = "  -3.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "612.41"
$ws.Range("E35").Value = "  -5.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.83"
$ws.Range("E36").Value = "  +8.40%  "
$ws.Range("E37").Value = "  -1.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.87"
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0483"
$ws.Range("E39").Value = "  +5.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "57.35"
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("E42").Value = "  +3.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.390.50"
$ws.Range("E43").Value = "  -0.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.322"
$ws.Range("E44").Value = "  -3.01%  "
$ws.Range("B45").Value = "PEPE"
$ws.Range("C45").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₃0715"
$ws.Range("E45").Value = "  +1.22%  "
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.00"
$ws.Range("E46").Value = "  +8.05%  "
$ws.Range("E47").Value = "  +0.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.62"
$ws.Range("E48").Value = "  +1.74%  "
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.02"
$ws.Range("E50").Value = "  -0.16%  "
